$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at C (old C "Feedback" column shifts to D)
$ws.Columns.Item(3).Insert()

# 2. Column widths: A (new) and C (new) get explicit widths; D keeps its old (shifted) width untouched
$ws.Columns.Item(1).ColumnWidth = 10.592447916666666
$ws.Columns.Item(3).ColumnWidth = 8.451822916666666

# 3. New row 5 data: B5 = 19, D5 = long feedback text, wrap text style, custom row height
#    (set this before "Over" so shared-string append order matches the target file)
$ws.Range("B5").Value = 19
$ws.Range("D5").Value = "task 1: 9.5/10`n- d) What about shape? explicitly mention position, color and shape in relationship to similarity and proximity`ntask 2: 9.5/10`n- initialisation of scatterplot does not work properly`n"
$ws.Range("D5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 90.75

# 4. New header cell C1 = "Over"
$ws.Range("C1").Value = "Over"

# 5. Row 13 total formula stays, but will recompute automatically once B5 is added
$ws.Range("B13").Formula = "=SUM(B2:B11)"

# 6. New row 14: "Over" total possible points
$ws.Range("A14").Value = "Over"
$ws.Range("B14").Formula = "=COUNT(B2:B11)*20"

# 7. New row 15: Percentage row with percent-formatted value and bonus note
$ws.Range("A15").Value = "Percentage"
$percentStyle = $wb.Styles.Add("Percent")
$ws.Range("B15").Style = "Percent"
$ws.Range("B15").Formula = "=B13/B14"
$ws.Range("B15").NumberFormat = "0%"
$ws.Range("D15").Value = "90% for 0.3 bonus"

# 8. Bold header / label cells (apply individually - Union/multi-area ranges are unreliable)
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("A13").Font.Bold = $true
$ws.Range("A14").Font.Bold = $true
$ws.Range("A15").Font.Bold = $true

# 9. Selection cursor position as in the diff
$ws.Range("D20").Select()
